$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(46).Insert()

$ws.Cells.Item(46, 1).Value = 4
$ws.Cells.Item(46, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(46, 3).Value = "Los Lagos"
$ws.Cells.Item(46, 4).Value = 44497
$ws.Cells.Item(46, 5).Value = 10
$ws.Cells.Item(46, 6).Value = 100112037
$ws.Cells.Item(46, 7).Value = "Cebollín"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 100
$ws.Cells.Item(46, 11).Value = 5000
$ws.Cells.Item(46, 12).Value = 6000
$ws.Cells.Item(46, 13).Value = 5500
$ws.Cells.Item(46, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(46, 15).Value = "Región Metropolitana"
$ws.Cells.Item(46, 16).Value = 153
$ws.Cells.Item(46, 17).Value = 36
$ws.Cells.Item(46, 18).Value = "Hortaliza"
